$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = '23.639.99'
$ws.Range("E2").Value2 = '  +1.90%  '

# Row 3
$ws.Range("D3").Value2 = '1.660.51'
$ws.Range("E3").Value2 = '  +3.01%  '

# Row 4
$ws.Range("D4").Value2 = '0.9977'
$ws.Range("E4").Value2 = '  -0.52%  '

# Row 5
$ws.Range("D5").Value2 = '0.9995'
$ws.Range("E5").Value2 = '  -0.26%  '

# Row 6
$ws.Range("D6").Value2 = '302.32'
$ws.Range("E6").Value2 = '  -0.03%  '

# Row 7
$ws.Range("D7").Value2 = '0.3837'
$ws.Range("E7").Value2 = '  +1.40%  '

# Row 8
$ws.Range("B8").Value2 = 'Cardano'
$ws.Range("C8").Value2 = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value2 = '0.3600'
$ws.Range("E8").Value2 = '  +2.05%  '

# Row 9
$ws.Range("B9").Value2 = 'OKB'
$ws.Range("C9").Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value2 = '51.08'
$ws.Range("E9").Value2 = '  -1.37%  '

# Row 10
$ws.Range("D10").Value2 = '1.245'
$ws.Range("E10").Value2 = '  +3.92%  '

# Row 11
$ws.Range("D11").Value2 = '0.08192'
$ws.Range("E11").Value2 = '  +1.32%  '

# Row 12
$ws.Range("D12").Value2 = '0.9975'
$ws.Range("E12").Value2 = '  -0.60%  '

# Row 13
$ws.Range("D13").Value2 = '22.46'
$ws.Range("E13").Value2 = '  +2.17%  '

# Row 14
$ws.Range("D14").Value2 = '6.510'
$ws.Range("E14").Value2 = '  +2.38%  '

# Row 15
$ws.Range("D15").Value2 = '7.514'
$ws.Range("E15").Value2 = '  +4.00%  '

# Row 16
$ws.Range("E16").Value2 = '  +1.68%  '

# Row 17
$ws.Range("D17").Value2 = '1.653.78'
$ws.Range("E17").Value2 = '  +2.58%  '

# Row 18
$ws.Range("D18").Value2 = '97.56'
$ws.Range("E18").Value2 = '  +3.57%  '

# Row 19
$ws.Range("D19").Value2 = '0.06977'
$ws.Range("E19").Value2 = '  +0.93%  '

# Row 20
$ws.Range("D20").Value2 = '6.838'
$ws.Range("E20").Value2 = '  +5.07%  '

# Row 21
$ws.Range("D21").Value2 = '17.73'
$ws.Range("E21").Value2 = '  +3.47%  '

# Row 22
$ws.Range("D22").Value2 = '0.9993'
$ws.Range("E22").Value2 = '  -0.21%  '

# Row 23
$ws.Range("D23").Value2 = '12.73'
$ws.Range("E23").Value2 = '  +3.38%  '

# Row 24
$ws.Range("D24").Value2 = '23.640.79'
$ws.Range("E24").Value2 = '  +1.95%  '

# Row 25
$ws.Range("D25").Value2 = '2.508'
$ws.Range("E25").Value2 = '  +0.15%  '

# Row 26
$ws.Range("D26").Value2 = '2.999'
$ws.Range("E26").Value2 = '  -0.49%  '

# Row 27
$ws.Range("E27").Value2 = '  +1.84%  '

# Row 28
$ws.Range("D28").Value2 = '151.77'
$ws.Range("E28").Value2 = '  +0.55%  '

# Row 29
$ws.Range("D29").Value2 = '5.236'
$ws.Range("E29").Value2 = '  +0.07%  '

# Row 30
$ws.Range("D30").Value2 = '134.01'
$ws.Range("E30").Value2 = '  +1.30%  '

# Row 31
$ws.Range("B31").Value2 = 'Filecoin'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value2 = '7.226'
$ws.Range("E31").Value2 = '  +11.84%  '

# Row 32
$ws.Range("B32").Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D32").Value2 = '1.839.01'
$ws.Range("E32").Value2 = '  +2.49%  '

# Row 33
$ws.Range("D33").Value2 = '2.242'
$ws.Range("E33").Value2 = '  +6.69%  '

# Row 34
$ws.Range("D34").Value2 = '12.01'
$ws.Range("E34").Value2 = '  +5.57%  '

# Row 35
$ws.Range("E35").Value2 = '  -1.11%  '

# Row 36
$ws.Range("D36").Value2 = '0.02807'
$ws.Range("E36").Value2 = '  +4.07%  '

# Row 37
$ws.Range("D37").Value2 = '6.142'
$ws.Range("E37").Value2 = '  +5.71%  '

# Row 38
$ws.Range("D38").Value2 = '0.2500'
$ws.Range("E38").Value2 = '  +2.35%  '

# Row 39
$ws.Range("D39").Value2 = '0.08794'
$ws.Range("E39").Value2 = '  +0.82%  '

# Row 40
$ws.Range("D40").Value2 = '0.07031'
$ws.Range("E40").Value2 = '  +1.53%  '

# Row 41
$ws.Range("D41").Value2 = '13.25'
$ws.Range("E41").Value2 = '  +10.92%  '

# Row 42
$ws.Range("D42").Value2 = '0.7044'
$ws.Range("E42").Value2 = '  +2.88%  '

# Row 43
$ws.Range("D43").Value2 = '1.337'
$ws.Range("E43").Value2 = '  +1.23%  '

# Row 44
$ws.Range("E44").Value2 = '  +5.42%  '

# Row 45
$ws.Range("D45").Value2 = '0.6558'
$ws.Range("E45").Value2 = '  +4.44%  '

# Row 46
$ws.Range("D46").Value2 = '0.9993'
$ws.Range("E46").Value2 = '  -0.23%  '

# Row 47
$ws.Range("D47").Value2 = '2.311'
$ws.Range("E47").Value2 = '  +3.22%  '

# Row 48
$ws.Range("D48").Value2 = '3.959'
$ws.Range("E48").Value2 = '  +0.46%  '

# Row 49
$ws.Range("D49").Value2 = '0.07950'
$ws.Range("E49").Value2 = '  +1.19%  '

# Row 50
$ws.Range("D50").Value2 = '128.02'
$ws.Range("E50").Value2 = '  +0.83%  '

# Row 51
$ws.Range("D51").Value2 = '1.198'
$ws.Range("E51").Value2 = '  +3.06%  '
